$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Build the full data block (rows 2-18, columns A-H) reflecting the scraper refresh
$data = New-Object 'object[,]' 17,8
$data[0,0] = "2026-02-09 19:20:46"
$data[0,1] = "初回 生成AIを業務に組み込むシステム開発"
$data[0,2] = "システム開発"
$data[0,3] = "300,000 円 ~ 500,000 円 / 固定"
$data[0,4] = "期限情報なし"
$data[0,5] = "https://www.lancers.jp/work/detail/5488757"
$data[0,6] = 410
$data[0,7] = "🔥AI,Ai ◆開発,システム開発"

$data[1,0] = "2026-02-09 19:20:46"
$data[1,1] = "【未来予測】パラレルワールドAIシステム開発の依頼"
$data[1,2] = "システム開発"
$data[1,3] = "100,000 円 ~ 200,000 円 / 固定"
$data[1,4] = "期限情報なし"
$data[1,5] = "https://www.lancers.jp/work/detail/5488301"
$data[1,6] = 403
$data[1,7] = "🔥AI,Ai ◆開発,システム開発"

$data[2,0] = "2026-02-09 19:20:46"
$data[2,1] = "製造業向け設備要件定義書の自動生成AIシステムの開発・DB設計支援エンジニア(AI/バックエンド)"
$data[2,2] = "システム開発"
$data[2,3] = "300,000 円 ~ 500,000 円 / 固定"
$data[2,4] = "期限情報なし"
$data[2,5] = "https://www.lancers.jp/work/detail/5473648"
$data[2,6] = 390
$data[2,7] = "🔥AI,Ai ◆開発"

$data[3,0] = "2026-02-09 19:20:46"
$data[3,1] = "競馬AIの開発ができる方、もしくはすでに開発済みの方"
$data[3,2] = "システム開発"
$data[3,3] = "300,000 円 ~ 500,000 円 / 固定"
$data[3,4] = "期限情報なし"
$data[3,5] = "https://www.lancers.jp/work/detail/5488810"
$data[3,6] = 375
$data[3,7] = "🔥AI,Ai ◆開発"

$data[4,0] = "2026-02-09 19:20:46"
$data[4,1] = "【急募】AIロボット・エージェント動作生成&販売プラットフォーム開発"
$data[4,2] = "システム開発"
$data[4,3] = "200,000 円 ~ 300,000 円 / 固定"
$data[4,4] = "期限情報なし"
$data[4,5] = "https://www.lancers.jp/work/detail/5488299"
$data[4,6] = 368
$data[4,7] = "🔥AI,Ai ◆開発"

$data[5,0] = "2026-02-09 19:20:46"
$data[5,1] = "【急募】パーソナルAI開発プロジェクトの依頼"
$data[5,2] = "システム開発"
$data[5,3] = "200,000 円 ~ 300,000 円 / 固定"
$data[5,4] = "期限情報なし"
$data[5,5] = "https://www.lancers.jp/work/detail/5488286"
$data[5,6] = 368
$data[5,7] = "🔥AI,Ai ◆開発"

$data[6,0] = "2026-02-09 19:20:46"
$data[6,1] = "【急募】AIシミュレーション相性チェックサービス開発者募集"
$data[6,2] = "システム開発"
$data[6,3] = "100,000 円 ~ 200,000 円 / 固定"
$data[6,4] = "期限情報なし"
$data[6,5] = "https://www.lancers.jp/work/detail/5488266"
$data[6,6] = 368
$data[6,7] = "🔥AI,Ai ◆開発"

$data[7,0] = "2026-02-09 19:20:46"
$data[7,1] = "【無在庫ツール開発】KeepaAPIとbaseAPIを活用したシステム構築"
$data[7,2] = "システム開発"
$data[7,3] = "1,000 ~ 5,000 円 / 固定"
$data[7,4] = "期限情報なし"
$data[7,5] = "https://www.lancers.jp/work/detail/5488392"
$data[7,6] = 320
$data[7,7] = "🔥API ◆ツール,開発"

$data[8,0] = "2026-02-09 19:20:46"
$data[8,1] = "※急募:Flutterによる業務アプリの開発(+next.js)"
$data[8,2] = "システム開発"
$data[8,3] = "300,000 円 ~ 500,000 円 / 固定"
$data[8,4] = "期限情報なし"
$data[8,5] = "https://www.lancers.jp/work/detail/5488271"
$data[8,6] = 225
$data[8,7] = "🔥Next.js ◆開発 ◇アプリ"

$data[9,0] = "2026-02-09 19:20:46"
$data[9,1] = "【Zapier保守・運用サポート】既存フローの管理・調整をお任せできる方募集(時給1,200円程度)"
$data[9,2] = "システム開発"
$data[9,3] = "50,000 円 ~ 100,000 円 / 固定"
$data[9,4] = "期限情報なし"
$data[9,5] = "https://www.lancers.jp/work/detail/5488168"
$data[9,6] = 213
$data[9,7] = "🔥API ◇管理"

$data[10,0] = "2026-02-09 19:20:46"
$data[10,1] = "【急募】多言語動画吹替・字幕一括生成システム開発"
$data[10,2] = "システム開発"
$data[10,3] = "100,000 円 ~ 200,000 円 / 固定"
$data[10,4] = "期限情報なし"
$data[10,5] = "https://www.lancers.jp/work/detail/5488276"
$data[10,6] = 118
$data[10,7] = "◆開発,システム開発"

$data[11,0] = "2026-02-09 19:20:46"
$data[11,1] = "【急募】クリニックの自動シフト調整システムをの開発お手伝いください!"
$data[11,2] = "システム開発"
$data[11,3] = "20,000 円 ~ 50,000 円 / 固定"
$data[11,4] = "期限情報なし"
$data[11,5] = "https://www.lancers.jp/work/detail/5488573"
$data[11,6] = 78
$data[11,7] = "◆開発"

$data[12,0] = "2026-02-09 19:20:46"
$data[12,1] = "製造業DXプロダクト開発のプロダクトマネージャー募集"
$data[12,2] = "システム開発"
$data[12,3] = "300,000 円 ~ 500,000 円 / 固定"
$data[12,4] = "期限情報なし"
$data[12,5] = "https://www.lancers.jp/work/detail/5468432"
$data[12,6] = 75
$data[12,7] = "◆開発"

$data[13,0] = "2026-02-09 19:20:46"
$data[13,1] = "スプレッドシート(Apps Script)で作業時間をボタン1つで計測・集計できる仕組みの開発"
$data[13,2] = "システム開発"
$data[13,3] = "50,000 円 ~ 100,000 円 / 固定"
$data[13,4] = "期限情報なし"
$data[13,5] = "https://www.lancers.jp/work/detail/5488743"
$data[13,6] = 68
$data[13,7] = "◆開発"

$data[14,0] = "2026-02-09 19:20:46"
$data[14,1] = "【急募】フルリモートでの自治体向け勤怠管理システム構築"
$data[14,2] = "システム開発"
$data[14,3] = "500,000 円 ~ 1,000,000 円 / 固定"
$data[14,4] = "期限情報なし"
$data[14,5] = "https://www.lancers.jp/work/detail/5488565"
$data[14,6] = 60
$data[14,7] = "◇管理"

$data[15,0] = "2026-02-09 19:20:46"
$data[15,1] = "養鰻管理Excelの判断ロジック(給餌)を理解し、継続的に伴走できる方を募集"
$data[15,2] = "システム開発"
$data[15,3] = "100,000 円 ~ 200,000 円 / 固定"
$data[15,4] = "期限情報なし"
$data[15,5] = "https://www.lancers.jp/work/detail/5488109"
$data[15,6] = 38
$data[15,7] = "◇管理"

$data[16,0] = "2026-02-09 19:20:46"
$data[16,1] = "【3月/フルリモート】証券バッチシステム改修案件(PL/SQL・Linux)"
$data[16,2] = "システム開発"
$data[16,3] = "500,000 円 ~ 1,000,000 円 / 固定"
$data[16,4] = "期限情報なし"
$data[16,5] = "https://www.lancers.jp/work/detail/5488543"
$data[16,6] = 40
$data[16,7] = $null

$ws.Range("A2:H18").Value = $data

# Rebuild hyperlinks on column F so each link points at the correct detail URL
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5488757", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488757") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5488301", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488301") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5473648", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5473648") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5488810", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488810") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5488299", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488299") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5488286", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488286") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5488266", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488266") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5488392", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488392") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5488271", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488271") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5488168", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488168") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5488276", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488276") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5488573", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488573") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5468432", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5468432") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5488743", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488743") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5488565", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488565") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5488109", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488109") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5488543", [System.Type]::Missing, [System.Type]::Missing, "https://www.lancers.jp/work/detail/5488543") | Out-Null

Write-Output "Updated rows 2-18 with refreshed scrape data and hyperlinks"
